# Updates market-price / profit figures on several sheets of the Kujata
# Profits workbook (scheduled runner refresh). Each block targets one
# worksheet + leve row, rewriting the currentAveragePrice*/LevePrice*/
# LeveProfit* columns (H:N) with refreshed marketboard data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3162.5
$ws.Range("I28").Value = 2884.7222
$ws.Range("J28").Value = 4412.5
$ws.Range("K28").Value = 2884.7222
$ws.Range("L28").Value = 4412.5
$ws.Range("M28").Value = -2399.7222
$ws.Range("N28").Value = -5382.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2457.3333
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2686
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2686
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3768

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6542537
$ws.Range("I132").Value = 7580043
$ws.Range("J132").Value = 21072.285
$ws.Range("K132").Value = 22740129
$ws.Range("L132").Value = 63216.855
$ws.Range("M132").Value = -22737599
$ws.Range("N132").Value = -68276.855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2100.81
$ws.Range("I138").Value = 952.0909
$ws.Range("J138").Value = 2424.8076
$ws.Range("K138").Value = 2856.2727
$ws.Range("L138").Value = 7274.4228
$ws.Range("M138").Value = 2283.7273
$ws.Range("N138").Value = -17554.4228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 111112930
$ws.Range("I61").Value = 166668060
$ws.Range("J61").Value = 2666
$ws.Range("K61").Value = 166668060
$ws.Range("L61").Value = 2666
$ws.Range("M61").Value = -166667848
$ws.Range("N61").Value = -3090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1370.16
$ws.Range("I74").Value = 935.8421
$ws.Range("K74").Value = 935.8421
$ws.Range("M74").Value = -61.84209999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 23333.334
$ws.Range("J76").Value = 23333.334
$ws.Range("L76").Value = 23333.334
$ws.Range("N76").Value = -24009.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1370.16
$ws.Range("I77").Value = 935.8421
$ws.Range("K77").Value = 4679.2105
$ws.Range("M77").Value = -311.2105000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 23333.334
$ws.Range("J79").Value = 23333.334
$ws.Range("L79").Value = 23333.334
$ws.Range("N79").Value = -25673.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2906.7805
$ws.Range("I132").Value = 2711.4375
$ws.Range("J132").Value = 3601.3333
$ws.Range("K132").Value = 8134.3125
$ws.Range("L132").Value = 10803.9999
$ws.Range("M132").Value = -5604.3125
$ws.Range("N132").Value = -15863.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 111112930
$ws.Range("I136").Value = 166668060
$ws.Range("J136").Value = 2666
$ws.Range("K136").Value = 500004180
$ws.Range("L136").Value = 7998
$ws.Range("M136").Value = -500001630
$ws.Range("N136").Value = -13098

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3730.4167
$ws.Range("I134").Value = 903.3077
$ws.Range("J134").Value = 11080.9
$ws.Range("K134").Value = 2709.9231
$ws.Range("L134").Value = 33242.7
$ws.Range("M134").Value = -174.9231
$ws.Range("N134").Value = -38312.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8199.556
$ws.Range("I58").Value = 1766.6666
$ws.Range("J58").Value = 11416
$ws.Range("K58").Value = 1766.6666
$ws.Range("L58").Value = 11416
$ws.Range("M58").Value = -1563.6666
$ws.Range("N58").Value = -11822

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 33544.668
$ws.Range("I108").Value = 30000
$ws.Range("K108").Value = 30000
$ws.Range("M108").Value = -26160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2520.7693
$ws.Range("I132").Value = 2916.5
$ws.Range("J132").Value = 1887.6
$ws.Range("K132").Value = 8749.5
$ws.Range("L132").Value = 5662.799999999999
$ws.Range("M132").Value = -6219.5
$ws.Range("N132").Value = -10722.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 18520110
$ws.Range("I134").Value = 1700
$ws.Range("J134").Value = 50001410
$ws.Range("K134").Value = 5100
$ws.Range("L134").Value = 150004230
$ws.Range("M134").Value = -2565
$ws.Range("N134").Value = -150009300

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8199.556
$ws.Range("I136").Value = 1766.6666
$ws.Range("J136").Value = 11416
$ws.Range("K136").Value = 5299.9998
$ws.Range("L136").Value = 34248
$ws.Range("M136").Value = -2749.9998
$ws.Range("N136").Value = -39348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 881300.9
$ws.Range("J141").Value = 881300.9
$ws.Range("L141").Value = 881300.9
$ws.Range("N141").Value = -891660.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3220.4
$ws.Range("I94").Value = 1001
$ws.Range("K94").Value = 3003
$ws.Range("M94").Value = -2327

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19260594
$ws.Range("J131").Value = 36764.57
$ws.Range("L131").Value = 110293.71
$ws.Range("N131").Value = -120373.71

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 769.25
$ws.Range("I97").Value = 646.7692
$ws.Range("K97").Value = 646.7692
$ws.Range("M97").Value = -150.7692

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1089.2222
$ws.Range("I122").Value = 1100.375
$ws.Range("K122").Value = 3301.125
$ws.Range("M122").Value = -851.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5160.154
$ws.Range("I132").Value = 5538.7144
$ws.Range("J132").Value = 4718.5
$ws.Range("K132").Value = 16616.1432
$ws.Range("L132").Value = 14155.5
$ws.Range("M132").Value = -14086.1432
$ws.Range("N132").Value = -19215.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1301.5
$ws.Range("I68").Value = 1301.5
$ws.Range("J68").Value = 1301.5
$ws.Range("K68").Value = 1301.5
$ws.Range("L68").Value = 1301.5
$ws.Range("M68").Value = -552.5
$ws.Range("N68").Value = -2799.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1301.5
$ws.Range("I71").Value = 1301.5
$ws.Range("J71").Value = 1301.5
$ws.Range("K71").Value = 6507.5
$ws.Range("L71").Value = 6507.5
$ws.Range("M71").Value = -2763.5
$ws.Range("N71").Value = -13995.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 11000
$ws.Range("J97").Value = 11000
$ws.Range("L97").Value = 11000
$ws.Range("N97").Value = -12982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1668.5714
$ws.Range("I136").Value = 1359.4117
$ws.Range("K136").Value = 4078.2351
$ws.Range("M136").Value = -1528.2351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2388.7568
$ws.Range("I132").Value = 2302.724
$ws.Range("J132").Value = 2700.625
$ws.Range("K132").Value = 6908.172
$ws.Range("L132").Value = 8101.875
$ws.Range("M132").Value = -4378.172
$ws.Range("N132").Value = -13161.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1905.625
$ws.Range("I136").Value = 1571.909
$ws.Range("K136").Value = 4715.727000000001
$ws.Range("M136").Value = -2165.727000000001

Write-Output "Kujata_Profits: refreshed 29 leve rows across 8 sheets"

